$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Transmitance (%) readings for rows 3-18 (column B), wavelengths
# 350.0 .. 1100.0 nm - values for the extra-ordinary ray (polarimetric
# acquisitions) channel.
$newTransmitance = @(
    39.60,  # 350.0 nm
    61.93,  # 400.0 nm
    75.66,  # 450.0 nm
    87.12,  # 500.0 nm
    88.20,  # 550.0 nm
    88.67,  # 600.0 nm
    87.36,  # 650.0 nm
    86.55,  # 700.0 nm
    76.48,  # 750.0 nm
    64.63,  # 800.0 nm
    48.95,  # 850.0 nm
    32.34,  # 900.0 nm
    17.37,  # 950.0 nm
    6.60,   # 1000.0 nm
    1.49,   # 1050.0 nm
    0.38    # 1100.0 nm
)

$startRow = 3
for ($i = 0; $i -lt $newTransmitance.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $newTransmitance[$i]
}

# The Transmitance column (B3:B18) drops its wrap-text formatting.
$ws.Range("B3:B18").WrapText = $false

# Move the active cell/selection to H7.
$ws.Range("H7").Select()
